$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (2021-07-10)
$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2818.026602407467

# Row 3 (2021-06-28)
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.260211312413533
